$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new (most recent) price record for "Bruselas (repollito)" needs to be added
# to the weekly log. The sheet keeps its rows in a fixed (non-chronological)
# order, and new observations are inserted near the top of the data block
# (row 23), pushing every subsequent row down by one. This naturally grows
# the used range from A1:R60 to A1:R61.

$ws.Rows.Item(23).Insert()

$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value = 44483
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 100112035
$ws.Range("G23").Value = "Bruselas (repollito)"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = 25000
$ws.Range("L23").Value = 25000
$ws.Range("M23").Value = 25000
$ws.Range("N23").Value = "$/malla 10 kilos"
$ws.Range("O23").Value = "Provincia de Quillota"
$ws.Range("P23").Value = 2500
$ws.Range("Q23").Value = 10
$ws.Range("R23").Value = "Hortaliza"
